# Edit script: rename "Requested quantity" headers, add "PO Forecast" sheet with forecast data.

$wb = $excel.ActiveWorkbook

# --- Step 1: rename the "Requested quantity" header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 2: add the new "PO Forecast" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match the page-margin defaults used by the other sheets in this workbook.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy header formatting (bold/border/center) from an existing header row.
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-column formatting (numFmt) down column A.
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A71").PasteSpecial(-4122)

# --- Step 3: write the header labels ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- Step 4: write the forecast data rows (A2:D71) ---
$data = New-Object 'object[,]' 70,4
$data[0,0]=44934.99999999999; $data[0,1]=88; $data[0,2]=-108.854384366215; $data[0,3]=275.158724460017
$data[1,0]=44948.99999999999; $data[1,1]=89; $data[1,2]=-126.8383184092235; $data[1,3]=298.8590288623307
$data[2,0]=44955.99999999999; $data[2,1]=89; $data[2,2]=-119.5274932426075; $data[2,3]=292.149823329933
$data[3,0]=44962.99999999999; $data[3,1]=89; $data[3,2]=-108.0314628133439; $data[3,3]=307.0159770603289
$data[4,0]=44969.99999999999; $data[4,1]=90; $data[4,2]=-125.3135781624699; $data[4,3]=284.5160825734419
$data[5,0]=44976.99999999999; $data[5,1]=90; $data[5,2]=-116.3761344615603; $data[5,3]=296.2207753587543
$data[6,0]=44983.99999999999; $data[6,1]=90; $data[6,2]=-104.8918406828645; $data[6,3]=288.7573207446924
$data[7,0]=44990.99999999999; $data[7,1]=91; $data[7,2]=-111.559370572663; $data[7,3]=292.7838451870339
$data[8,0]=44997.99999999999; $data[8,1]=91; $data[8,2]=-118.8041880718632; $data[8,3]=290.7129604213933
$data[9,0]=45004.99999999999; $data[9,1]=91; $data[9,2]=-94.43593675825663; $data[9,3]=282.2199849627162
$data[10,0]=45011.99999999999; $data[10,1]=91; $data[10,2]=-112.7378652609599; $data[10,3]=304.5893761376823
$data[11,0]=45018.99999999999; $data[11,1]=92; $data[11,2]=-114.4093094699469; $data[11,3]=282.5669847089095
$data[12,0]=45025.99999999999; $data[12,1]=92; $data[12,2]=-113.6491419244195; $data[12,3]=289.7708687530668
$data[13,0]=45039.99999999999; $data[13,1]=93; $data[13,2]=-97.01139669117465; $data[13,3]=303.5722708422591
$data[14,0]=45046.99999999999; $data[14,1]=93; $data[14,2]=-107.6226717494796; $data[14,3]=296.4200480272291
$data[15,0]=45053.99999999999; $data[15,1]=93; $data[15,2]=-99.04463850565307; $data[15,3]=294.5137073426836
$data[16,0]=45060.99999999999; $data[16,1]=93; $data[16,2]=-102.5061084767717; $data[16,3]=302.0125276218778
$data[17,0]=45088.99999999999; $data[17,1]=95; $data[17,2]=-122.1319343747721; $data[17,3]=295.6452328960368
$data[18,0]=45200.99999999999; $data[18,1]=99; $data[18,2]=-105.1762192632675; $data[18,3]=304.1624347966392
$data[19,0]=45207.99999999999; $data[19,1]=99; $data[19,2]=-116.0482406455967; $data[19,3]=287.3734145208845
$data[20,0]=45214.99999999999; $data[20,1]=100; $data[20,2]=-97.80162382344376; $data[20,3]=299.7948468174988
$data[21,0]=45221.99999999999; $data[21,1]=100; $data[21,2]=-84.10472048292473; $data[21,3]=301.2326826281622
$data[22,0]=45228.99999999999; $data[22,1]=100; $data[22,2]=-105.7147101948732; $data[22,3]=308.6584920863075
$data[23,0]=45235.99999999999; $data[23,1]=101; $data[23,2]=-90.48346492406849; $data[23,3]=290.7942378808796
$data[24,0]=45242.99999999999; $data[24,1]=101; $data[24,2]=-110.4618744798705; $data[24,3]=307.6436974148938
$data[25,0]=45249.99999999999; $data[25,1]=101; $data[25,2]=-102.9780539884564; $data[25,3]=307.1680729146943
$data[26,0]=45256.99999999999; $data[26,1]=101; $data[26,2]=-100.1308411819537; $data[26,3]=313.3121944005274
$data[27,0]=45270.99999999999; $data[27,1]=102; $data[27,2]=-113.6604332210568; $data[27,3]=312.0105862469396
$data[28,0]=45277.99999999999; $data[28,1]=102; $data[28,2]=-96.45201112897823; $data[28,3]=320.6695344885768
$data[29,0]=45298.99999999999; $data[29,1]=103; $data[29,2]=-109.6593780346819; $data[29,3]=308.792746506144
$data[30,0]=45305.99999999999; $data[30,1]=104; $data[30,2]=-87.18749090593032; $data[30,3]=311.0814477874171
$data[31,0]=45312.99999999999; $data[31,1]=104; $data[31,2]=-88.8158078305181; $data[31,3]=309.0012386040615
$data[32,0]=45319.99999999999; $data[32,1]=104; $data[32,2]=-97.69867173010562; $data[32,3]=307.1468937233942
$data[33,0]=45326.99999999999; $data[33,1]=104; $data[33,2]=-93.04088249012713; $data[33,3]=296.245746636367
$data[34,0]=45333.99999999999; $data[34,1]=105; $data[34,2]=-96.8261954004381; $data[34,3]=308.1367003736223
$data[35,0]=45340.99999999999; $data[35,1]=105; $data[35,2]=-82.67537729951502; $data[35,3]=296.2742170248105
$data[36,0]=45347.99999999999; $data[36,1]=105; $data[36,2]=-96.22375632538527; $data[36,3]=307.1467673351411
$data[37,0]=45354.99999999999; $data[37,1]=106; $data[37,2]=-103.9100252843713; $data[37,3]=290.0487180299617
$data[38,0]=45361.99999999999; $data[38,1]=106; $data[38,2]=-103.9708994642887; $data[38,3]=302.2787189152932
$data[39,0]=45368.99999999999; $data[39,1]=106; $data[39,2]=-109.3590055192621; $data[39,3]=297.0555990913369
$data[40,0]=45375.99999999999; $data[40,1]=106; $data[40,2]=-100.1329941143637; $data[40,3]=307.330510319949
$data[41,0]=45382.99999999999; $data[41,1]=107; $data[41,2]=-74.06080032506016; $data[41,3]=303.593119707712
$data[42,0]=45389.99999999999; $data[42,1]=107; $data[42,2]=-86.47391190333643; $data[42,3]=312.5376513486394
$data[43,0]=45417.99999999999; $data[43,1]=108; $data[43,2]=-88.63569785483213; $data[43,3]=306.0079595945595
$data[44,0]=45424.99999999999; $data[44,1]=108; $data[44,2]=-94.56156082360894; $data[44,3]=311.3779998316827
$data[45,0]=45431.99999999999; $data[45,1]=109; $data[45,2]=-90.4557822543066; $data[45,3]=310.4940767244863
$data[46,0]=45438.99999999999; $data[46,1]=109; $data[46,2]=-90.31737091916084; $data[46,3]=310.2767833594755
$data[47,0]=45445.99999999999; $data[47,1]=109; $data[47,2]=-85.58015799701258; $data[47,3]=311.2427440866518
$data[48,0]=45452.99999999999; $data[48,1]=110; $data[48,2]=-100.2149932818786; $data[48,3]=308.7312418541626
$data[49,0]=45459.99999999999; $data[49,1]=110; $data[49,2]=-89.62058055454821; $data[49,3]=301.7114465214304
$data[50,0]=45473.99999999999; $data[50,1]=110; $data[50,2]=-94.75450512373129; $data[50,3]=315.2928091759426
$data[51,0]=45480.99999999999; $data[51,1]=111; $data[51,2]=-87.43504238910758; $data[51,3]=303.1263945634133
$data[52,0]=45487.99999999999; $data[52,1]=111; $data[52,2]=-88.81930111746867; $data[52,3]=307.3438271682543
$data[53,0]=45536.99999999999; $data[53,1]=113; $data[53,2]=-79.91343339239349; $data[53,3]=307.3898479841085
$data[54,0]=45543.99999999999; $data[54,1]=113; $data[54,2]=-73.88044116081414; $data[54,3]=304.7512337227516
$data[55,0]=45550.99999999999; $data[55,1]=114; $data[55,2]=-96.40498904144303; $data[55,3]=302.1411548017934
$data[56,0]=45557.99999999999; $data[56,1]=114; $data[56,2]=-83.31050854683465; $data[56,3]=319.9491299094082
$data[57,0]=45599.99999999999; $data[57,1]=116; $data[57,2]=-70.12833767374167; $data[57,3]=310.3372444704871
$data[58,0]=45613.99999999999; $data[58,1]=116; $data[58,2]=-79.42392853470928; $data[58,3]=321.3186333868268
$data[59,0]=45620.99999999999; $data[59,1]=116; $data[59,2]=-76.04533790483293; $data[59,3]=313.9049100455761
$data[60,0]=45627.99999999999; $data[60,1]=117; $data[60,2]=-82.65843914353965; $data[60,3]=299.5654878459691
$data[61,0]=45634.99999999999; $data[61,1]=117; $data[61,2]=-78.39495532860994; $data[61,3]=296.6157190831212
$data[62,0]=45641.99999999999; $data[62,1]=117; $data[62,2]=-78.04211832853589; $data[62,3]=312.9644880628748
$data[63,0]=45648.99999999999; $data[63,1]=118; $data[63,2]=-73.75780134972034; $data[63,3]=320.313600161386
$data[64,0]=45655.99999999999; $data[64,1]=118; $data[64,2]=-93.5627307487963; $data[64,3]=310.3142448918629
$data[65,0]=45662.99999999999; $data[65,1]=118; $data[65,2]=-87.44487249057968; $data[65,3]=310.6789883092803
$data[66,0]=45669.99999999999; $data[66,1]=118; $data[66,2]=-78.08095793969851; $data[66,3]=313.5177349946
$data[67,0]=45676.99999999999; $data[67,1]=119; $data[67,2]=-88.09279110839657; $data[67,3]=303.2772170948106
$data[68,0]=45683.99999999999; $data[68,1]=119; $data[68,2]=-74.62066444073531; $data[68,3]=322.5252757776443
$data[69,0]=45690.99999999999; $data[69,1]=119; $data[69,2]=-91.00075628962355; $data[69,3]=312.5160697176315
$newSheet.Range("A2:D71").Value = $data

Write-Output "PO Forecast sheet added with $($newSheet.UsedRange.Rows.Count) rows"
